$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set values for row 5 (Дороган Михаела): C5, D5, E5 = 5
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 5

# Set values for row 12 (Кравчук Мария): C12, D12, E12 = 5
$ws.Range("C12").Value = 5
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 5

# Set value for row 18 (Петров Иван): C18 = 5
$ws.Range("C18").Value = 5

# Update the active selection to C18
$ws.Range("C18").Select()
